$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Correct the value of F1 (was erroneously "c5", should be "F1" like the
# rest of the header/label column F2:F6)
$ws.Range("F1").Value = "F1"

# Update the active selection on Sheet1 to reflect the new cursor position
$ws.Range("I22").Select()
